$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-06-16 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-17 Saturday", 2) | Out-Null
$d.Content.Find.Execute("20×22=440", $true, $false, $false, $false, $false, $true, 1, $false, "59×13=767", 2) | Out-Null
$d.Content.Find.Execute("49×54=2646", $true, $false, $false, $false, $false, $true, 1, $false, "60×85=5100", 2) | Out-Null
$d.Content.Find.Execute("86×93=7998", $true, $false, $false, $false, $false, $true, 1, $false, "99×62=6138", 2) | Out-Null
$d.Content.Find.Execute("42×15=630", $true, $false, $false, $false, $false, $true, 1, $false, "49×71=3479", 2) | Out-Null
$d.Content.Find.Execute("51×98=4998", $true, $false, $false, $false, $false, $true, 1, $false, "38×30=1140", 2) | Out-Null
$d.Content.Find.Execute("74×78=5772", $true, $false, $false, $false, $false, $true, 1, $false, "65×59=3835", 2) | Out-Null
$d.Content.Find.Execute("59×42=2478", $true, $false, $false, $false, $false, $true, 1, $false, "63×25=1575", 2) | Out-Null
$d.Content.Find.Execute("86×33=2838", $true, $false, $false, $false, $false, $true, 1, $false, "39×28=1092", 2) | Out-Null
$d.Content.Find.Execute("13×35=455", $true, $false, $false, $false, $false, $true, 1, $false, "10×77=770", 2) | Out-Null
$d.Content.Find.Execute("29×56=1624", $true, $false, $false, $false, $false, $true, 1, $false, "27×66=1782", 2) | Out-Null
$d.Content.Find.Execute("85×68=5780", $true, $false, $false, $false, $false, $true, 1, $false, "36×29=1044", 2) | Out-Null
$d.Content.Find.Execute("57×35=1995", $true, $false, $false, $false, $false, $true, 1, $false, "80×70=5600", 2) | Out-Null
$d.Content.Find.Execute("58×42=2436", $true, $false, $false, $false, $false, $true, 1, $false, "28×53=1484", 2) | Out-Null
$d.Content.Find.Execute("62×83=5146", $true, $false, $false, $false, $false, $true, 1, $false, "88×19=1672", 2) | Out-Null
$d.Content.Find.Execute("74×87=6438", $true, $false, $false, $false, $false, $true, 1, $false, "75×71=5325", 2) | Out-Null
$d.Content.Find.Execute("21×62=1302", $true, $false, $false, $false, $false, $true, 1, $false, "59×26=1534", 2) | Out-Null
$d.Content.Find.Execute("35×50=1750", $true, $false, $false, $false, $false, $true, 1, $false, "56×90=5040", 2) | Out-Null
$d.Content.Find.Execute("73×24=1752", $true, $false, $false, $false, $false, $true, 1, $false, "84×47=3948", 2) | Out-Null
$d.Content.Find.Execute("83×77=6391", $true, $false, $false, $false, $false, $true, 1, $false, "24×26=624", 2) | Out-Null
$d.Content.Find.Execute("33×27=891", $true, $false, $false, $false, $false, $true, 1, $false, "43×92=3956", 2) | Out-Null
$d.Content.Find.Execute("60×84=5040", $true, $false, $false, $false, $false, $true, 1, $false, "55×90=4950", 2) | Out-Null
$d.Content.Find.Execute("20×51=1020", $true, $false, $false, $false, $false, $true, 1, $false, "85×44=3740", 2) | Out-Null
$d.Content.Find.Execute("48×24=1152", $true, $false, $false, $false, $false, $true, 1, $false, "71×50=3550", 2) | Out-Null
$d.Content.Find.Execute("75×50=3750", $true, $false, $false, $false, $false, $true, 1, $false, "47×99=4653", 2) | Out-Null
$d.Content.Find.Execute("25×19=475", $true, $false, $false, $false, $false, $true, 1, $false, "28×62=1736", 2) | Out-Null
$d.Content.Find.Execute("98×23=2254", $true, $false, $false, $false, $false, $true, 1, $false, "44×86=3784", 2) | Out-Null
$d.Content.Find.Execute("71×98=6958", $true, $false, $false, $false, $false, $true, 1, $false, "45×30=1350", 2) | Out-Null
$d.Content.Find.Execute("78×12=936", $true, $false, $false, $false, $false, $true, 1, $false, "17×96=1632", 2) | Out-Null
$d.Content.Find.Execute("66×16=1056", $true, $false, $false, $false, $false, $true, 1, $false, "75×25=1875", 2) | Out-Null
$d.Content.Find.Execute("79×26=2054", $true, $false, $false, $false, $false, $true, 1, $false, "26×15=390", 2) | Out-Null
$d.Content.Find.Execute("75×47=3525", $true, $false, $false, $false, $false, $true, 1, $false, "97×95=9215", 2) | Out-Null
$d.Content.Find.Execute("55×36=1980", $true, $false, $false, $false, $false, $true, 1, $false, "24×31=744", 2) | Out-Null
$d.Content.Find.Execute("54×31=1674", $true, $false, $false, $false, $false, $true, 1, $false, "53×68=3604", 2) | Out-Null
$d.Content.Find.Execute("92×79=7268", $true, $false, $false, $false, $false, $true, 1, $false, "54×53=2862", 2) | Out-Null
$d.Content.Find.Execute("32×95=3040", $true, $false, $false, $false, $false, $true, 1, $false, "66×22=1452", 2) | Out-Null
$d.Content.Find.Execute("61×53=3233", $true, $false, $false, $false, $false, $true, 1, $false, "40×100=4000", 2) | Out-Null
$d.Content.Find.Execute("44×91=4004", $true, $false, $false, $false, $false, $true, 1, $false, "50×65=3250", 2) | Out-Null
$d.Content.Find.Execute("66×68=4488", $true, $false, $false, $false, $false, $true, 1, $false, "72×84=6048", 2) | Out-Null
$d.Content.Find.Execute("67×83=5561", $true, $false, $false, $false, $false, $true, 1, $false, "97×36=3492", 2) | Out-Null
$d.Content.Find.Execute("21×23=483", $true, $false, $false, $false, $false, $true, 1, $false, "68×12=816", 2) | Out-Null
$d.Content.Find.Execute("72×74=5328", $true, $false, $false, $false, $false, $true, 1, $false, "39×47=1833", 2) | Out-Null
$d.Content.Find.Execute("87×71=6177", $true, $false, $false, $false, $false, $true, 1, $false, "62×22=1364", 2) | Out-Null
$d.Content.Find.Execute("94×52=4888", $true, $false, $false, $false, $false, $true, 1, $false, "18×10=180", 2) | Out-Null
$d.Content.Find.Execute("67×44=2948", $true, $false, $false, $false, $false, $true, 1, $false, "89×84=7476", 2) | Out-Null
$d.Content.Find.Execute("57×92=5244", $true, $false, $false, $false, $false, $true, 1, $false, "65×92=5980", 2) | Out-Null
$d.Content.Find.Execute("66×38=2508", $true, $false, $false, $false, $false, $true, 1, $false, "95×95=9025", 2) | Out-Null
$d.Content.Find.Execute("100×12=1200", $true, $false, $false, $false, $false, $true, 1, $false, "89×70=6230", 2) | Out-Null
$d.Content.Find.Execute("15×89=1335", $true, $false, $false, $false, $false, $true, 1, $false, "38×83=3154", 2) | Out-Null
$d.Content.Find.Execute("41×10=410", $true, $false, $false, $false, $false, $true, 1, $false, "90×26=2340", 2) | Out-Null
$d.Content.Find.Execute("45×43=1935", $true, $false, $false, $false, $false, $true, 1, $false, "23×53=1219", 2) | Out-Null
$d.Content.Find.Execute("58×56=3248", $true, $false, $false, $false, $false, $true, 1, $false, "86×62=5332", 2) | Out-Null
$d.Content.Find.Execute("99×70=6930", $true, $false, $false, $false, $false, $true, 1, $false, "63×75=4725", 2) | Out-Null
$d.Content.Find.Execute("75×38=2850", $true, $false, $false, $false, $false, $true, 1, $false, "80×51=4080", 2) | Out-Null
$d.Content.Find.Execute("88×96=8448", $true, $false, $false, $false, $false, $true, 1, $false, "47×59=2773", 2) | Out-Null
$d.Content.Find.Execute("11×94=1034", $true, $false, $false, $false, $false, $true, 1, $false, "58×66=3828", 2) | Out-Null
$d.Content.Find.Execute("83×50=4150", $true, $false, $false, $false, $false, $true, 1, $false, "80×65=5200", 2) | Out-Null
$d.Content.Find.Execute("90×16=1440", $true, $false, $false, $false, $false, $true, 1, $false, "82×99=8118", 2) | Out-Null
$d.Content.Find.Execute("97×19=1843", $true, $false, $false, $false, $false, $true, 1, $false, "84×12=1008", 2) | Out-Null
$d.Content.Find.Execute("19×52=988", $true, $false, $false, $false, $false, $true, 1, $false, "37×82=3034", 2) | Out-Null
$d.Content.Find.Execute("13×98=1274", $true, $false, $false, $false, $false, $true, 1, $false, "22×53=1166", 2) | Out-Null
$d.Content.Find.Execute("60×52=3120", $true, $false, $false, $false, $false, $true, 1, $false, "70×99=6930", 2) | Out-Null
$d.Content.Find.Execute("20×97=1940", $true, $false, $false, $false, $false, $true, 1, $false, "53×10=530", 2) | Out-Null
$d.Content.Find.Execute("40×28=1120", $true, $false, $false, $false, $false, $true, 1, $false, "44×78=3432", 2) | Out-Null
$d.Content.Find.Execute("73×71=5183", $true, $false, $false, $false, $false, $true, 1, $false, "84×50=4200", 2) | Out-Null
$d.Content.Find.Execute("46×22=1012", $true, $false, $false, $false, $false, $true, 1, $false, "49×98=4802", 2) | Out-Null
$d.Content.Find.Execute("98×93=9114", $true, $false, $false, $false, $false, $true, 1, $false, "81×15=1215", 2) | Out-Null
$d.Content.Find.Execute("58×22=1276", $true, $false, $false, $false, $false, $true, 1, $false, "97×62=6014", 2) | Out-Null
$d.Content.Find.Execute("58×93=5394", $true, $false, $false, $false, $false, $true, 1, $false, "23×63=1449", 2) | Out-Null
$d.Content.Find.Execute("16×84=1344", $true, $false, $false, $false, $false, $true, 1, $false, "32×99=3168", 2) | Out-Null
$d.Content.Find.Execute("73×90=6570", $true, $false, $false, $false, $false, $true, 1, $false, "85×76=6460", 2) | Out-Null
$d.Content.Find.Execute("55×30=1650", $true, $false, $false, $false, $false, $true, 1, $false, "78×50=3900", 2) | Out-Null
$d.Content.Find.Execute("42×63=2646", $true, $false, $false, $false, $false, $true, 1, $false, "93×46=4278", 2) | Out-Null
$d.Content.Find.Execute("11×29=319", $true, $false, $false, $false, $false, $true, 1, $false, "26×90=2340", 2) | Out-Null
$d.Content.Find.Execute("69×30=2070", $true, $false, $false, $false, $false, $true, 1, $false, "25×100=2500", 2) | Out-Null
$d.Content.Find.Execute("81×39=3159", $true, $false, $false, $false, $false, $true, 1, $false, "54×35=1890", 2) | Out-Null
$d.Content.Find.Execute("63×50=3150", $true, $false, $false, $false, $false, $true, 1, $false, "77×90=6930", 2) | Out-Null
$d.Content.Find.Execute("18×100=1800", $true, $false, $false, $false, $false, $true, 1, $false, "43×77=3311", 2) | Out-Null
$d.Content.Find.Execute("100×32=3200", $true, $false, $false, $false, $false, $true, 1, $false, "74×97=7178", 2) | Out-Null
$d.Content.Find.Execute("85×15=1275", $true, $false, $false, $false, $false, $true, 1, $false, "24×67=1608", 2) | Out-Null
$d.Content.Find.Execute("85×41=3485", $true, $false, $false, $false, $false, $true, 1, $false, "92×57=5244", 2) | Out-Null
$d.Content.Find.Execute("32×20=640", $true, $false, $false, $false, $false, $true, 1, $false, "37×77=2849", 2) | Out-Null
$d.Content.Find.Execute("97×63=6111", $true, $false, $false, $false, $false, $true, 1, $false, "63×95=5985", 2) | Out-Null
$d.Content.Find.Execute("55×62=3410", $true, $false, $false, $false, $false, $true, 1, $false, "67×25=1675", 2) | Out-Null
$d.Content.Find.Execute("51×12=612", $true, $false, $false, $false, $false, $true, 1, $false, "10×14=140", 2) | Out-Null
$d.Content.Find.Execute("11×64=704", $true, $false, $false, $false, $false, $true, 1, $false, "69×63=4347", 2) | Out-Null
$d.Content.Find.Execute("95×55=5225", $true, $false, $false, $false, $false, $true, 1, $false, "98×22=2156", 2) | Out-Null
$d.Content.Find.Execute("58×83=4814", $true, $false, $false, $false, $false, $true, 1, $false, "91×78=7098", 2) | Out-Null
$d.Content.Find.Execute("51×20=1020", $true, $false, $false, $false, $false, $true, 1, $false, "100×86=8600", 2) | Out-Null
$d.Content.Find.Execute("92×32=2944", $true, $false, $false, $false, $false, $true, 1, $false, "15×48=720", 2) | Out-Null
$d.Content.Find.Execute("25×21=525", $true, $false, $false, $false, $false, $true, 1, $false, "96×41=3936", 2) | Out-Null
$d.Content.Find.Execute("82×29=2378", $true, $false, $false, $false, $false, $true, 1, $false, "73×66=4818", 2) | Out-Null
$d.Content.Find.Execute("39×39=1521", $true, $false, $false, $false, $false, $true, 1, $false, "100×36=3600", 2) | Out-Null
$d.Content.Find.Execute("48×53=2544", $true, $false, $false, $false, $false, $true, 1, $false, "45×74=3330", 2) | Out-Null
$d.Content.Find.Execute("96×75=7200", $true, $false, $false, $false, $false, $true, 1, $false, "89×56=4984", 2) | Out-Null
$d.Content.Find.Execute("70×88=6160", $true, $false, $false, $false, $false, $true, 1, $false, "81×45=3645", 2) | Out-Null
$d.Content.Find.Execute("37×97=3589", $true, $false, $false, $false, $false, $true, 1, $false, "76×99=7524", 2) | Out-Null
$d.Content.Find.Execute("18×53=954", $true, $false, $false, $false, $false, $true, 1, $false, "100×70=7000", 2) | Out-Null
$d.Content.Find.Execute("11×16=176", $true, $false, $false, $false, $false, $true, 1, $false, "32×77=2464", 2) | Out-Null
$d.Content.Find.Execute("61×41=2501", $true, $false, $false, $false, $false, $true, 1, $false, "17×10=170", 2) | Out-Null
$d.Content.Find.Execute("32×80=2560", $true, $false, $false, $false, $false, $true, 1, $false, "26×64=1664", 2) | Out-Null
